$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rename the existing "birth circumstance" sheet to "birth distribution"
# ---------------------------------------------------------------------------
$wsBirthDist = $wb.Worksheets.Item(7)
$wsBirthDist.Name = "birth distribution"

# ---------------------------------------------------------------------------
# 2. Add the three new sheets, in order, right after "birth distribution"
# ---------------------------------------------------------------------------
$wsTimeBetween = $wb.Worksheets.Add($null, $wsBirthDist)
$wsTimeBetween.Name = "time between births"

$wsRRType = $wb.Worksheets.Add($null, $wsTimeBetween)
$wsRRType.Name = "RR birth by type"

$wsRRTime = $wb.Worksheets.Add($null, $wsRRType)
$wsRRTime.Name = "RR birth by time"

# ---------------------------------------------------------------------------
# 3. Populate "birth distribution" (order of writes controls shared-string
#    indices, so the order below matches the canonical sheet layout exactly)
# ---------------------------------------------------------------------------
$wsBirthDist.Range("A1").Value = "Type"
$wsBirthDist.Range("B1").Value = "<18 years"
$wsBirthDist.Range("C1").Value = "18-34 years"
$wsBirthDist.Range("D1").Value = "35-49 years"

$wsBirthDist.Range("A2").Value = "first"
$wsBirthDist.Range("B2").Value = 0.0543
$wsBirthDist.Range("C2").Value = 0.1711
$wsBirthDist.Range("D2").Value = 0.0003

$wsBirthDist.Range("A3").Value = "second or third"
$wsBirthDist.Range("B3").Value = 0.009
$wsBirthDist.Range("C3").Value = 0.3607
$wsBirthDist.Range("D3").Value = 0.0085

$wsBirthDist.Range("A4").Value = "greater than third"
$wsBirthDist.Range("B4").Value = 0.0
$wsBirthDist.Range("C4").Value = 0.2908
$wsBirthDist.Range("D4").Value = 0.1048

# ---------------------------------------------------------------------------
# 4. Populate "time between births"
# ---------------------------------------------------------------------------
$wsTimeBetween.Range("A1").Value = "first"
$wsTimeBetween.Range("B1").Value = "<18 months"
$wsTimeBetween.Range("C1").Value = "18-23 months"
$wsTimeBetween.Range("D1").Value = "<24 months"

$wsTimeBetween.Range("A2").Value = 0.2258
$wsTimeBetween.Range("B2").Value = 0.0705
$wsTimeBetween.Range("C2").Value = 0.134
$wsTimeBetween.Range("D2").Value = 0.5698

# ---------------------------------------------------------------------------
# 5. Populate "RR birth by type"
# ---------------------------------------------------------------------------
$wsRRType.Range("A1").Value = "Outcome"
$wsRRType.Range("B1").Value = "Type"
$wsRRType.Range("C1").Value = "<18 years"
$wsRRType.Range("D1").Value = "18-34 years"
$wsRRType.Range("E1").Value = "35-49 years"

$wsRRType.Range("A2").Value = "pretermSGA"
$wsRRType.Range("B2").Value = "first"
$wsRRType.Range("C2").Value = 3.14
$wsRRType.Range("D2").Value = 1.73
$wsRRType.Range("E2").Value = 1.73

$wsRRType.Range("B3").Value = "second or third"
$wsRRType.Range("C3").Value = 1.6
$wsRRType.Range("D3").Value = 1.0
$wsRRType.Range("E3").Value = 1.57

$wsRRType.Range("B4").Value = "greater than third"
$wsRRType.Range("C4").Value = 1.6
$wsRRType.Range("D4").Value = 1.0
$wsRRType.Range("E4").Value = 1.57

$wsRRType.Range("A5").Value = "pretermAGA"
$wsRRType.Range("B5").Value = "first"
$wsRRType.Range("C5").Value = 1.75
$wsRRType.Range("D5").Value = 1.75
$wsRRType.Range("E5").Value = 1.75

$wsRRType.Range("B6").Value = "second or third"
$wsRRType.Range("C6").Value = 1.4
$wsRRType.Range("D6").Value = 1.0
$wsRRType.Range("E6").Value = 1.33

$wsRRType.Range("B7").Value = "greater than third"
$wsRRType.Range("C7").Value = 1.4
$wsRRType.Range("D7").Value = 1.0
$wsRRType.Range("E7").Value = 1.33

$wsRRType.Range("A8").Value = "termSGA"
$wsRRType.Range("B8").Value = "first"
$wsRRType.Range("C8").Value = 1.52
$wsRRType.Range("D8").Value = 1.52
$wsRRType.Range("E8").Value = 1.52

$wsRRType.Range("B9").Value = "second or third"
$wsRRType.Range("C9").Value = 1.2
$wsRRType.Range("D9").Value = 1.0
$wsRRType.Range("E9").Value = 1.0

$wsRRType.Range("B10").Value = "greater than third"
$wsRRType.Range("C10").Value = 1.2
$wsRRType.Range("D10").Value = 1.0
$wsRRType.Range("E10").Value = 1.0

# ---------------------------------------------------------------------------
# 6. Populate "RR birth by time"
# ---------------------------------------------------------------------------
$wsRRTime.Range("A1").Value = "Outcome"
$wsRRTime.Range("B1").Value = "first"
$wsRRTime.Range("C1").Value = "<18 months"
$wsRRTime.Range("D1").Value = "18-23 months"
$wsRRTime.Range("E1").Value = "<24 months"

$wsRRTime.Range("A2").Value = "pretermSGA"
$wsRRTime.Range("B2").Value = 1.0
$wsRRTime.Range("C2").Value = 3.03
$wsRRTime.Range("D2").Value = 1.77
$wsRRTime.Range("E2").Value = 1.0

$wsRRTime.Range("A3").Value = "pretermAGA"
$wsRRTime.Range("B3").Value = 1.0
$wsRRTime.Range("C3").Value = 1.49
$wsRRTime.Range("D3").Value = 1.1
$wsRRTime.Range("E3").Value = 1.0

$wsRRTime.Range("A4").Value = "termSGA"
$wsRRTime.Range("B4").Value = 1.0
$wsRRTime.Range("C4").Value = 1.41
$wsRRTime.Range("D4").Value = 1.18
$wsRRTime.Range("E4").Value = 1.0
